$d = $word.ActiveDocument

# --- Change 1: the two runs "'zero-to-hero.html" + "' " are merged into a
# single run "'zero-to-hero.html' ". This is a purely cosmetic run-split
# change - the rendered text is identical, so no text-level edit is needed.

# --- Change 2: remove the paragraph "If you reach this point stop for a
# while!" together with the blank paragraph that follows it.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*If you reach this point stop for a while!*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $nextPara = $targetPara.Next()
    $startRange = $targetPara.Range.Start
    $endRange = $nextPara.Range.End
    $rangeToDelete = $d.Range($startRange, $endRange)
    $rangeToDelete.Delete()
}

# --- Change 3: "Overwrite the doGet method." -> "Override the doGet() method."
$find = $d.Content.Find
$find.ClearFormatting()
$null = $find.Execute("Overwrite the doGet method.", $false, $false, $false, $false, $false, $true, 1, $false, "Override the doGet() method.", 2)

# --- Change 4: "overwritten" -> "overridden" (both occurrences).
# Replace the first occurrence alone so we can relocate the hidden
# "_GoBack" bookmark right after it (matching where Word leaves it after
# the last text edit), then replace the remaining occurrence.
$find = $d.Content.Find
$find.ClearFormatting()
$null = $find.Execute("overwritten", $false, $false, $false, $false, $false, $true, 1, $false, "overridden", 1)

$firstOccurrenceEnd = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*overridden method build*") {
        $text = $p.Range.Text
        $idx = $text.IndexOf("overridden")
        $firstOccurrenceEnd = $p.Range.Start + $idx + "overridden".Length
        break
    }
}

$find = $d.Content.Find
$find.ClearFormatting()
$null = $find.Execute("overwritten", $false, $false, $false, $false, $false, $true, 1, $false, "overridden", 2)

# --- Move the hidden "_GoBack" bookmark to sit right after "overridden" in
# "In the overridden method build an html table ..." (its position after
# the author's final text edit in this area).
if ($firstOccurrenceEnd -ne $null) {
    $newRange = $d.Range($firstOccurrenceEnd, $firstOccurrenceEnd)
    $d.Bookmarks.Add("_GoBack", $newRange)
}
